$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '96.188.17'
$ws.Range("E2").Value = '  -0.89%  '

# Row 3
$ws.Range("D3").Value = '3.623.06'
$ws.Range("E3").Value = '  -1.99%  '

# Row 4
$ws.Range("D4").Value = '''2.77'
$ws.Range("E4").Value = '  +29.66%  '

# Row 5
$ws.Range("E5").Value = '  +0.02%  '

# Row 6
$ws.Range("D6").Value = '''224.54'
$ws.Range("E6").Value = '  -4.68%  '

# Row 7
$ws.Range("D7").Value = '''640.90'
$ws.Range("E7").Value = '  -2.25%  '

# Row 8
$ws.Range("D8").Value = '''0.426'
$ws.Range("E8").Value = '  -0.69%  '

# Row 9
$ws.Range("E9").Value = '  +12.87%  '

# Row 10
$ws.Range("E10").Value = '  +0.00%  '

# Row 11
$ws.Range("D11").Value = '3.616.79'
$ws.Range("E11").Value = '  -2.12%  '

# Row 12
$ws.Range("D12").Value = '''48.54'
$ws.Range("E12").Value = '  +8.48%  '

# Row 13
$ws.Range("D13").Value = '''0.216'
$ws.Range("E13").Value = '  +4.34%  '

# Row 14
$ws.Range("D14").Value = '''0.0000292'
$ws.Range("E14").Value = '  -4.11%  '

# Row 15
$ws.Range("D15").Value = '''6.55'
$ws.Range("E15").Value = '  -3.72%  '

# Row 16
$ws.Range("D16").Value = '4.297.96'
$ws.Range("E16").Value = '  -2.02%  '

# Row 17
$ws.Range("D17").Value = '95.842.53'
$ws.Range("E17").Value = '  -1.01%  '

# Row 18
$ws.Range("D18").Value = '''24.22'
$ws.Range("E18").Value = '  +30.51%  '

# Row 19
$ws.Range("D19").Value = '''9.00'
$ws.Range("E19").Value = '  -0.73%  '

# Row 20
$ws.Range("D20").Value = '''13.92'
$ws.Range("E20").Value = '  +7.62%  '

# Row 21
$ws.Range("D21").Value = '3.629.36'
$ws.Range("E21").Value = '  -1.84%  '

# Row 22
$ws.Range("D22").Value = '''0.291'
$ws.Range("E22").Value = '  +45.84%  '

# Row 23
$ws.Range("D23").Value = '''0.543'
$ws.Range("E23").Value = '  +0.97%  '

# Row 24
$ws.Range("B24").Value = 'Litecoin'
$ws.Range("C24").Value = 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc'
$ws.Range("D24").Value = '''135.37'
$ws.Range("E24").Value = '  +23.19%  '

# Row 25
$ws.Range("B25").Value = 'BitcoinCash'
$ws.Range("C25").Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
$ws.Range("D25").Value = '''524.30'
$ws.Range("E25").Value = '  +0.67%  '

# Row 26
$ws.Range("D26").Value = '''3.26'
$ws.Range("E26").Value = '  -4.86%  '

# Row 27
$ws.Range("D27").Value = '''0.0000204'
$ws.Range("E27").Value = '  -7.07%  '

# Row 28
$ws.Range("D28").Value = '''6.83'
$ws.Range("E28").Value = '  -0.53%  '

# Row 29
$ws.Range("D29").Value = '3.787.74'
$ws.Range("E29").Value = '  -2.67%  '

# Row 30
$ws.Range("D30").Value = '''12.90'
$ws.Range("E30").Value = '  -3.64%  '

# Row 31
$ws.Range("D31").Value = '''13.33'
$ws.Range("E31").Value = '  +6.61%  '

# Row 32
$ws.Range("D32").Value = '''3.11'
$ws.Range("E32").Value = '  +3.64%  '

# Row 33
$ws.Range("E33").Value = '  +0.22%  '

# Row 34
$ws.Range("D34").Value = '''0.633'
$ws.Range("E34").Value = '  +7.05%  '

# Row 35
$ws.Range("D35").Value = '''0.184'
$ws.Range("E35").Value = '  -1.86%  '

# Row 36
$ws.Range("D36").Value = '''33.26'
$ws.Range("E36").Value = '  +1.86%  '

# Row 37
$ws.Range("E37").Value = '  +0.24%  '

# Row 38
$ws.Range("D38").Value = '''1.79'
$ws.Range("E38").Value = '  -1.15%  '

# Row 39
$ws.Range("B39").Value = 'Filecoin'
$ws.Range("C39").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D39").Value = '''7.41'
$ws.Range("E39").Value = '  +9.40%  '

# Row 40
$ws.Range("D40").Value = '''0.540'
$ws.Range("E40").Value = '  +10.67%  '

# Row 41
$ws.Range("B41").Value = 'USDe'
$ws.Range("C41").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range("D41").Value = '''1.00'
$ws.Range("E41").Value = '  -0.03%  '

# Row 42
$ws.Range("D42").Value = '''592.18'
$ws.Range("E42").Value = '  -6.01%  '

# Row 43
$ws.Range("D43").Value = '''8.36'
$ws.Range("E43").Value = '  -4.08%  '

# Row 44
$ws.Range("D44").Value = '''0.0531'
$ws.Range("E44").Value = '  +18.37%  '

# Row 45
$ws.Range("D45").Value = '''41.73'
$ws.Range("E45").Value = '  +5.43%  '

# Row 46
$ws.Range("B46").Value = 'ARBITRUM'
$ws.Range("C46").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D46").Value = '''0.977'
$ws.Range("E46").Value = '  +2.49%  '

# Row 47
$ws.Range("B47").Value = 'Kaspa'
$ws.Range("C47").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D47").Value = '''0.158'
$ws.Range("E47").Value = '  -5.25%  '

# Row 48
$ws.Range("D48").Value = '''1.97'
$ws.Range("E48").Value = '  -2.17%  '

# Row 49
$ws.Range("B49").Value = 'Aave'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D49").Value = '''238.31'
$ws.Range("E49").Value = '  +16.77%  '

# Row 50
$ws.Range("B50").Value = 'Cosmos'
$ws.Range("C50").Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Range("D50").Value = '''9.25'
$ws.Range("E50").Value = '  +6.84%  '

# Row 51
$ws.Range("D51").Value = '''2.30'
$ws.Range("E51").Value = '  -3.55%  '
